$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.379.76"
$ws.Cells.Item(2, 5).Value = "  +0.69%  "
$ws.Cells.Item(3, 4).Value = "3.194.25"
$ws.Cells.Item(3, 5).Value = "  +0.19%  "
$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = "  -0.15%  "
$ws.Cells.Item(5, 4).Value = "'606.36"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +2.17%  "
$ws.Cells.Item(6, 4).Value = "'154.99"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +0.73%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -0.07%  "
$ws.Cells.Item(8, 4).Value = "3.198.62"
$ws.Cells.Item(8, 5).Value = "  +0.34%  "
$ws.Cells.Item(9, 4).Value = "'0.545"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +2.23%  "
$ws.Cells.Item(10, 4).Value = "'0.160"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +0.08%  "
$ws.Cells.Item(11, 4).Value = "'5.64"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -7.35%  "
$ws.Cells.Item(12, 4).Value = "'0.509"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  -0.46%  "
$ws.Cells.Item(13, 4).Value = "'0.0000268"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -0.69%  "
$ws.Cells.Item(14, 4).Value = "'38.59"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -0.66%  "
$ws.Cells.Item(15, 4).Value = "3.714.78"
$ws.Cells.Item(15, 5).Value = "  +0.16%  "
$ws.Cells.Item(16, 4).Value = "66.334.11"
$ws.Cells.Item(16, 5).Value = "  +0.66%  "
$ws.Cells.Item(17, 4).Value = "'7.40"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +0.00%  "
$ws.Cells.Item(18, 4).Value = "3.191.79"
$ws.Cells.Item(18, 5).Value = "  +0.00%  "
$ws.Cells.Item(19, 4).Value = "'0.113"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +1.16%  "
$ws.Cells.Item(20, 4).Value = "'510.29"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +0.05%  "
$ws.Cells.Item(21, 4).Value = "'15.34"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +0.21%  "
$ws.Cells.Item(22, 4).Value = "'0.733"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -0.75%  "
$ws.Cells.Item(23, 4).Value = "'8.07"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +1.24%  "
$ws.Cells.Item(24, 4).Value = "'14.87"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -2.13%  "
$ws.Cells.Item(25, 4).Value = "'84.61"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -0.20%  "
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  +0.12%  "
$ws.Cells.Item(27, 4).Value = "'3.00"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +0.71%  "
$ws.Cells.Item(28, 4).Value = "'9.22"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -1.01%  "
$ws.Cells.Item(29, 4).Value = "'2.41"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +6.29%  "
$ws.Cells.Item(30, 4).Value = "'3.06"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +6.24%  "
$ws.Cells.Item(31, 4).Value = "'7.05"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +2.11%  "
$ws.Cells.Item(32, 4).Value = "'28.09"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -0.51%  "
$ws.Cells.Item(33, 4).Value = "'1.00"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  +0.10%  "
$ws.Cells.Item(34, 4).Value = "'1.18"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -3.60%  "
$ws.Cells.Item(35, 4).Value = "'6.51"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -0.48%  "
$ws.Cells.Item(36, 4).Value = "'510.78"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +5.90%  "
$ws.Cells.Item(37, 4).Value = "'55.77"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +1.85%  "
$ws.Cells.Item(38, 4).Value = "'0.0928"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  +2.88%  "
$ws.Cells.Item(39, 4).Value = "0.0₃0738"
$ws.Cells.Item(39, 5).Value = "  +13.90%  "
$ws.Cells.Item(40, 4).Value = "'0.0418"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +0.47%  "
$ws.Cells.Item(41, 4).Value = "'0.128"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +4.98%  "
$ws.Cells.Item(42, 4).Value = "'2.96"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +1.50%  "
$ws.Cells.Item(43, 4).Value = "'8.78"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.36%  "
$ws.Cells.Item(44, 4).Value = "'0.300"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +1.32%  "
$ws.Cells.Item(45, 4).Value = "'2.49"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +3.05%  "
$ws.Cells.Item(46, 4).Value = "2.837.87"
$ws.Cells.Item(46, 5).Value = "  -3.10%  "
$ws.Cells.Item(47, 4).Value = "'28.22"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -0.60%  "
$ws.Cells.Item(48, 4).Value = "'2.39"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +2.77%  "
$ws.Cells.Item(49, 4).Value = "'0.999"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -0.11%  "
$ws.Cells.Item(50, 4).Value = "'0.117"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +0.72%  "
$ws.Cells.Item(51, 4).Value = "'2.62"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -0.06%  "
